# Insert a new data row at row 10 (pushing existing rows 10..126 down to 11..127)
# and populate it with the new record described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = 2
$ws.Range("B10").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C10").Value = "Coquimbo"
$ws.Range("D10").Value = (Get-Date -Year 2023 -Month 4 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 100112030
$ws.Range("G10").Value = "Poroto granado"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 700
$ws.Range("K10").Value = 28000
$ws.Range("L10").Value = 30000
$ws.Range("M10").Value = 29000
$ws.Range("N10").Value = "`$/malla 25 kilos"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 1160
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"
